$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.692.62'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.805.35'
$ws.Range('E3').Value = '  -2.43%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'230.94"
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'38.94"
$ws.Range('E8').Value = '  -7.08%  '
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').Value = "'0.0989"
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('D12').Value = '2.065.75'
$ws.Range('E12').Value = '  -2.51%  '
$ws.Range('D13').Value = '1.809.67'
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').Value = "'0.659"
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = "'10.84"
$ws.Range('E15').Value = '  -5.05%  '
$ws.Range('D16').Value = "'4.53"
$ws.Range('E16').Value = '  -3.12%  '
$ws.Range('D17').Value = '34.682.88'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').Value = "'239.03"
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('E21').Value = '  -2.97%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('D25').Value = "'172.30"
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('D26').Value = "'7.67"
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('E28').Value = '  -2.63%  '
$ws.Range('E29').Value = '  +9.10%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = "'3.97"
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('D32').Value = "'0.0541"
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').Value = "'3.92"
$ws.Range('E33').Value = '  -3.15%  '
$ws.Range('D34').Value = "'1.24"
$ws.Range('E34').Value = '  +13.89%  '
$ws.Range('E35').Value = '  -5.66%  '
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('D37').Value = "'90.65"
$ws.Range('E37').Value = '  -7.62%  '
$ws.Range('E38').Value = '  +4.81%  '
$ws.Range('D39').Value = '1.306.09'
$ws.Range('E39').Value = '  -3.59%  '
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('D41').Value = "'2.46"
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').Value = "'0.953"
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('D43').Value = "'14.16"
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('E44').Value = '  -10.16%  '
$ws.Range('E45').Value = '  -5.39%  '
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('D48').Value = '1.993.70'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'1.01"
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0670"
$ws.Range('E50').Value = '  +7.85%  '
$ws.Range('D51').Value = "'98.45"
$ws.Range('E51').Value = '  -4.95%  '
